$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add the four new "product kit" rows under the existing quick-order rows ---
# Row 11: a lone kit SKU line (no quantity / unit yet)
$ws.Range("A11").Value = "product-kit-1"

# Row 12: kit SKU with a quantity
$ws.Range("A12").Value = "product-kit-1"
$ws.Range("B12").Value = 1

# Row 13: kit SKU with quantity + unit
$ws.Range("A13").Value = "product-kit-1"
$ws.Range("B13").Value = 2
$ws.Range("C13").Value = "milliliter"

# Row 14: second kit SKU with quantity + unit
$ws.Range("A14").Value = "product-kit-2"
$ws.Range("B14").Value = 1
$ws.Range("C14").Value = "milliliter"

# Rows keep the sheet's natural (auto) row height; nudge it to match the
# recalculated auto-height for the new rows.
$ws.Rows("11").RowHeight = 12.8
$ws.Rows("12").RowHeight = 12.8
$ws.Rows("13").RowHeight = 12.8
$ws.Rows("14").RowHeight = 12.8

# --- Selection moves to the newly entered quantity cell ---
$ws.Range("B13").Select() | Out-Null

# --- Minor page-margin nudge (header/footer) ---
$ps = $ws.PageSetup
$ps.HeaderMargin = 36.850393700787386
$ps.FooterMargin = 36.850393700787386
